# Insert a new weekly data row at row 283 (pushing existing rows 283..302
# down to 284..303) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(283).Insert()

$ws.Range("A283").Value = 6
$ws.Range("B283").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C283").Value = "Metropolitana"
$ws.Range("D283").Value = 45021
$ws.Range("E283").Value = 13
$ws.Range("F283").Value = 100112001
$ws.Range("G283").Value = "Berenjena"
$ws.Range("H283").Value = "Sin especificar"
$ws.Range("I283").Value = "Primera"
$ws.Range("J283").Value = 500
$ws.Range("K283").Value = 6000
$ws.Range("L283").Value = 7000
$ws.Range("M283").Value = 6540
$ws.Range("N283").Value = "`$/caja 50 unidades"
$ws.Range("O283").Value = "Región de Arica y Parinacota"
$ws.Range("P283").Value = 131
$ws.Range("Q283").Value = 50
$ws.Range("R283").Value = "Hortaliza"

$ws.Range("D283").NumberFormat = "YYYY-MM-DD HH:MM:SS"
